$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.060.40"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "'3.317.02"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'551.30"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").Value = "'172.02"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("E7").Value = "  +1.14%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'3.308.34"
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("D10").Value = "'0.169"
$ws.Range("E10").Value = "  +6.42%  "
$ws.Range("D11").Value = "'0.628"
$ws.Range("E11").Value = "  +1.76%  "
$ws.Range("D12").Value = "'53.36"
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("D13").Value = "'0.0000279"
$ws.Range("E13").Value = "  +3.86%  "
$ws.Range("D14").Value = "'9.02"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").Value = "'3.849.39"
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("E16").Value = "  +2.38%  "
$ws.Range("D17").Value = "'18.06"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").Value = "'3.323.82"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("D19").Value = "'63.952.86"
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("D20").Value = "'11.65"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").Value = "'0.982"
$ws.Range("E21").Value = "  +1.79%  "
$ws.Range("D22").Value = "'451.66"
$ws.Range("E22").Value = "  +6.65%  "
$ws.Range("E23").Value = "  +8.52%  "
$ws.Range("D24").Value = "'4.06"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").Value = "'86.99"
$ws.Range("E25").Value = "  +3.81%  "
$ws.Range("D26").Value = "'13.71"
$ws.Range("E26").Value = "  +4.40%  "
$ws.Range("E27").Value = "  +1.69%  "
$ws.Range("D28").Value = "'10.66"
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("D29").Value = "'8.55"
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("D30").Value = "'30.75"
$ws.Range("E30").Value = "  +4.59%  "
$ws.Range("D31").Value = "'6.52"
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("D33").Value = "'563.76"
$ws.Range("E33").Value = "  -4.51%  "
$ws.Range("D34").Value = "'60.67"
$ws.Range("E34").Value = "  +4.46%  "
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E37").Value = "  -2.04%  "
$ws.Range("D38").Value = "'3.50"
$ws.Range("E38").Value = "  +2.74%  "
$ws.Range("D39").Value = "'35.09"
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("D40").Value = "'0.365"
$ws.Range("E40").Value = "  +0.66%  "
$ws.Range("D41").Value = "'0.0₃0726"
$ws.Range("E41").Value = "  -1.90%  "
$ws.Range("D42").Value = "'3.055.15"
$ws.Range("E42").Value = "  -0.85%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0412"
$ws.Range("E43").Value = "  +2.14%  "
$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").Value = "'2.74"
$ws.Range("E44").Value = "  -1.81%  "
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("E46").Value = "  +3.08%  "
$ws.Range("D47").Value = "'2.42"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("D49").Value = "'141.73"
$ws.Range("E49").Value = "  +6.85%  "
$ws.Range("E50").Value = "  -2.69%  "
$ws.Range("D51").Value = "'8.10"
$ws.Range("E51").Value = "  +0.48%  "
